# Apply the diff: insert two new leading columns (data, loja), shift
# existing columns right, replace row 2 with new product data, and
# remove the old row 3 entirely.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two columns at the front so existing A:G data shifts to C:I.
$ws.Range("A:B").Insert()

# Delete the (now stale) third data row - was row 3, still row 3 after
# the column insert since rows are unaffected by a column insert.
$ws.Range("3:3").Delete()

# --- Header row (row 1) ---
$ws.Cells.Item(1, 1).Value = "data"
$ws.Cells.Item(1, 2).Value = "loja"

# Match the header formatting (bold, border, centered) used by the rest
# of row 1 by copying the format from column C onto the new A/B cells.
$ws.Range("C1").Copy()
$ws.Range("A1:B1").PasteSpecial(-4122)

# --- Data row (row 2) ---
$ws.Cells.Item(2, 1).Value = "30/07/2024"
$ws.Cells.Item(2, 2).Value = "ecomonline11"
$ws.Cells.Item(2, 3).Value = "Controle Longa Distancia Jfa K1200 Alcance De 1200 Metros"
$ws.Cells.Item(2, 4).Value = "K1200"
$ws.Cells.Item(2, 5).Value = 56.16
$ws.Cells.Item(2, 6).Value = "Baixo"
$ws.Cells.Item(2, 7).Value = "NA"
$ws.Cells.Item(2, 8).Value = "classico"
$ws.Cells.Item(2, 9).Value = "https://www.mercadolivre.com.br/controle-longa-distancia-jfa-k1200-alcance-de-1200-metros/p/MLB34245679?pdp_filters=seller_id:1568119549#searchVariation=MLB34245679&position=1&search_layout=stack&type=product&tracking_id=6c16a85a-68fd-437c-a9c0-84cd32267bbd"
